$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.753.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.451.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.19"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.37"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.555"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.449.86"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.887.80"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.715.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.451.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.19"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.33"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0738"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.01"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.75%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.46"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.784"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "269.94"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.97"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.584"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.27"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0908"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0484"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0210"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.73"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.723.96"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.77%  "
